$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "368.74") are not auto-converted into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "50.922.49"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "2.889.14"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "368.74"
$ws.Range("E5").Value = "  +3.91%  "
$ws.Range("D6").Value = "101.69"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("E7").Value = "  -4.98%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("D10").Value = "36.54"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "0.0831"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").Value = "3.347.68"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "7.31"
$ws.Range("E15").Value = "  -5.14%  "
$ws.Range("D16").Value = "2.880.73"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "0.927"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("D18").Value = "50.912.84"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "3.25"
$ws.Range("E19").Value = "  -6.08%  "
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -4.97%  "
$ws.Range("D21").Value = "12.75"
$ws.Range("E21").Value = "  -6.61%  "
$ws.Range("D22").Value = "0.0₃0939"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "67.76"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").Value = "257.87"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("D29").Value = "6.92"
$ws.Range("E29").Value = "  -9.54%  "
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").Value = "9.81"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").Value = "6.00"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").Value = "34.27"
$ws.Range("E34").Value = "  -6.47%  "
$ws.Range("D35").Value = "50.48"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "0.0413"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").Value = "3.02"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").Value = "16.84"
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("E41").Value = "  -8.20%  "
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("D43").Value = "21.79"
$ws.Range("E43").Value = "  -6.17%  "
$ws.Range("D44").Value = "117.94"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").Value = "2.020.06"
$ws.Range("E46").Value = "  -4.80%  "
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  -6.55%  "
$ws.Range("E48").Value = "  -6.98%  "
$ws.Range("D49").Value = "3.178.04"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").Value = "0.0308"
$ws.Range("E51").Value = "  -11.16%  "
